# msz - mandatory fields checks part 1
#
# Adds a new "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields"
# test-case row (columns A-D) to the process sheet, then mirrors the
# column-width auto-fit and cell-selection side effects Excel performs
# after typing the new, wider strings into row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row of data -------------------------------------------------
$ws.Range("A3").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields"
$ws.Range("B3").Value = "Vehicle Page check for open mandatory fields"
$ws.Range("C3").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields_FillMake"
$ws.Range("D3").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields_CheckFilledMake"

# --- column widths: the new, longer strings in A/C/D (and the wider
# content now visible in G) push those columns' auto-fit width out ----
$ws.Columns.Item(1).ColumnWidth = 55.666666666666664
$ws.Columns.Item(3).ColumnWidth = 73.0
$ws.Columns.Item(4).ColumnWidth = 70.66666666666667
$ws.Columns.Item(7).ColumnWidth = 26.0

# --- leave the selection on the newly-filled cell ---------------------
$ws.Range("B3").Select() | Out-Null

Write-Output "applied mandatory-fields row + column width updates"
